$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - sheet1
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 394
$wsExhibition.Range("F4").Value = 0

# Sheet "全部类型" (All types) - sheet4
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 0
$wsAll.Range("F5").Value = 19
$wsAll.Range("F7").Value = 0
$wsAll.Range("F8").Value = 146
$wsAll.Range("F9").Value = 63
